$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 2 (the "com.singleton.strechy / taxi game / redvelvetmichael@gmail.com" review),
# shifting all subsequent rows up by one.
$ws.Rows.Item(2).Delete()

$ws.Range("A2").Select()
